$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D29").Value = 44845
$ws.Range("K29").Value = 35000
$ws.Range("L29").Value = 35000
$ws.Range("M29").Value = 35000
$ws.Range("O29").Value = 'Perú'
$ws.Range("P29").Value = 1400
$ws.Range("D30").Value = 44624
$ws.Range("J30").Value = 40
$ws.Range("K30").Value = 30000
$ws.Range("L30").Value = 30000
$ws.Range("M30").Value = 30000
$ws.Range("N30").Value = '$/saco 25 kilos'
$ws.Range("O30").Value = 'Región Metropolitana'
$ws.Range("P30").Value = 1200
$ws.Range("D31").Value = 45027
$ws.Range("K31").Value = 25000
$ws.Range("L31").Value = 25000
$ws.Range("M31").Value = 25000
$ws.Range("P31").Value = 1000
$ws.Range("D32").Value = 44873
$ws.Range("K32").Value = 45000
$ws.Range("L32").Value = 45000
$ws.Range("M32").Value = 45000
$ws.Range("N32").Value = '$/malla 25 kilos'
$ws.Range("O32").Value = 'Perú'
$ws.Range("P32").Value = 1800
$ws.Range("D33").Value = 44642
$ws.Range("J33").Value = 60
$ws.Range("K33").Value = 36000
$ws.Range("L33").Value = 36000
$ws.Range("M33").Value = 36000
$ws.Range("N33").Value = '$/saco 25 kilos'
$ws.Range("O33").Value = 'Región Metropolitana'
$ws.Range("P33").Value = 1440
$ws.Range("D34").Value = 44383
$ws.Range("J34").Value = 40
$ws.Range("K34").Value = 32000
$ws.Range("L34").Value = 32000
$ws.Range("M34").Value = 32000
$ws.Range("N34").Value = '$/malla 25 kilos'
$ws.Range("O34").Value = 'Perú'
$ws.Range("P34").Value = 1280
$ws.Range("D35").Value = 44218
$ws.Range("J35").Value = 30
$ws.Range("K35").Value = 40000
$ws.Range("L35").Value = 40000
$ws.Range("M35").Value = 40000
$ws.Range("N35").Value = '$/saco 25 kilos'
$ws.Range("O35").Value = 'Región Metropolitana'
$ws.Range("P35").Value = 1600
$ws.Range("D36").Value = 44327
$ws.Range("K36").Value = 45000
$ws.Range("L36").Value = 45000
$ws.Range("M36").Value = 45000
$ws.Range("P36").Value = 1800
$ws.Range("D37").Value = 44260
$ws.Range("K37").Value = 30000
$ws.Range("L37").Value = 30000
$ws.Range("M37").Value = 30000
$ws.Range("P37").Value = 1200
$ws.Range("D38").Value = 44722
$ws.Range("J38").Value = 40
$ws.Range("K38").Value = 27000
$ws.Range("L38").Value = 27000
$ws.Range("M38").Value = 27000
$ws.Range("N38").Value = '$/malla 25 kilos'
$ws.Range("O38").Value = 'Perú'
$ws.Range("P38").Value = 1080
$ws.Range("D39").Value = 44358
$ws.Range("J39").Value = 35
$ws.Range("K39").Value = 35000
$ws.Range("L39").Value = 35000
$ws.Range("M39").Value = 35000
$ws.Range("O39").Value = 'Provincia de Limarí'
$ws.Range("P39").Value = 1400
$ws.Range("D40").Value = 44649
$ws.Range("J40").Value = 50
$ws.Range("K40").Value = 30000
$ws.Range("L40").Value = 30000
$ws.Range("M40").Value = 30000
$ws.Range("N40").Value = '$/saco 25 kilos'
$ws.Range("O40").Value = 'Región Metropolitana'
$ws.Range("P40").Value = 1200
$ws.Range("D41").Value = 44761
$ws.Range("K41").Value = 37000
$ws.Range("L41").Value = 37000
$ws.Range("M41").Value = 37000
$ws.Range("N41").Value = '$/malla 25 kilos'
$ws.Range("O41").Value = 'Perú'
$ws.Range("P41").Value = 1480
$ws.Range("D42").Value = 45006
$ws.Range("J42").Value = 45
$ws.Range("K42").Value = 32000
$ws.Range("L42").Value = 32000
$ws.Range("M42").Value = 32000
$ws.Range("N42").Value = '$/saco 25 kilos'
$ws.Range("O42").Value = 'Región Metropolitana'
$ws.Range("P42").Value = 1280
$ws.Range("D43").Value = 44915
$ws.Range("J43").Value = 35
$ws.Range("K43").Value = 35000
$ws.Range("L43").Value = 35000
$ws.Range("M43").Value = 35000
$ws.Range("P43").Value = 1400
$ws.Range("D44").Value = 44442
$ws.Range("K44").Value = 47000
$ws.Range("L44").Value = 47000
$ws.Range("M44").Value = 47000
$ws.Range("N44").Value = '$/malla 25 kilos'
$ws.Range("O44").Value = 'Perú'
$ws.Range("P44").Value = 1880
$ws.Range("D45").Value = 45009
$ws.Range("J45").Value = 40
$ws.Range("K45").Value = 30000
$ws.Range("L45").Value = 30000
$ws.Range("M45").Value = 30000
$ws.Range("O45").Value = 'Región Metropolitana'
$ws.Range("P45").Value = 1200
$ws.Range("D46").Value = 44285
$ws.Range("J46").Value = 30
$ws.Range("K46").Value = 40000
$ws.Range("L46").Value = 40000
$ws.Range("M46").Value = 40000
$ws.Range("N46").Value = '$/saco 25 kilos'
$ws.Range("P46").Value = 1600
$ws.Range("D47").Value = 44698
$ws.Range("J47").Value = 45
$ws.Range("K47").Value = 27000
$ws.Range("L47").Value = 27000
$ws.Range("M47").Value = 27000
$ws.Range("N47").Value = '$/malla 25 kilos'
$ws.Range("O47").Value = 'Perú'
$ws.Range("P47").Value = 1080
$ws.Range("D48").Value = 44236
$ws.Range("J48").Value = 30
$ws.Range("K48").Value = 32000
$ws.Range("L48").Value = 32000
$ws.Range("M48").Value = 32000
$ws.Range("N48").Value = '$/saco 25 kilos'
$ws.Range("O48").Value = 'Región Metropolitana'
$ws.Range("P48").Value = 1280
$ws.Range("D49").Value = 44859
$ws.Range("J49").Value = 35
$ws.Range("K49").Value = 36000
$ws.Range("L49").Value = 36000
$ws.Range("M49").Value = 36000
$ws.Range("N49").Value = '$/malla 25 kilos'
$ws.Range("O49").Value = 'Perú'
$ws.Range("P49").Value = 1440
$ws.Range("D50").Value = 44299
$ws.Range("J50").Value = 30
$ws.Range("K50").Value = 38000
$ws.Range("L50").Value = 38000
$ws.Range("M50").Value = 38000
$ws.Range("N50").Value = '$/saco 25 kilos'
$ws.Range("O50").Value = 'Región Metropolitana'
$ws.Range("P50").Value = 1520
$ws.Range("D51").Value = 44295
$ws.Range("K51").Value = 30000
$ws.Range("L51").Value = 30000
$ws.Range("M51").Value = 30000
$ws.Range("P51").Value = 1200
$ws.Range("D52").Value = 44400
$ws.Range("J52").Value = 40
$ws.Range("N52").Value = '$/malla 25 kilos'
$ws.Range("O52").Value = 'Perú'
$ws.Range("D53").Value = 44376
$ws.Range("D54").Value = 44691
$ws.Range("N54").Value = '$/saco 25 kilos'
$ws.Range("O54").Value = 'Región Metropolitana'
$ws.Range("D55").Value = 44628
$ws.Range("J55").Value = 50
$ws.Range("K55").Value = 38000
$ws.Range("L55").Value = 38000
$ws.Range("M55").Value = 38000
$ws.Range("P55").Value = 1520
$ws.Range("D56").Value = 44309
$ws.Range("J56").Value = 30
$ws.Range("K56").Value = 40000
$ws.Range("L56").Value = 40000
$ws.Range("M56").Value = 40000
$ws.Range("P56").Value = 1600
$ws.Range("D57").Value = 44985
$ws.Range("J57").Value = 45
$ws.Range("K57").Value = 30000
$ws.Range("L57").Value = 30000
$ws.Range("M57").Value = 30000
$ws.Range("P57").Value = 1200
$ws.Range("D58").Value = 45002
$ws.Range("J58").Value = 40
$ws.Range("K58").Value = 33000
$ws.Range("L58").Value = 33000
$ws.Range("M58").Value = 33000
$ws.Range("P58").Value = 1320
$ws.Range("D59").Value = 44278
$ws.Range("J59").Value = 30
$ws.Range("K59").Value = 36000
$ws.Range("L59").Value = 36000
$ws.Range("M59").Value = 36000
$ws.Range("P59").Value = 1440
$ws.Range("D60").Value = 45013
$ws.Range("J60").Value = 40
$ws.Range("K60").Value = 30000
$ws.Range("L60").Value = 30000
$ws.Range("M60").Value = 30000
$ws.Range("P60").Value = 1200
$ws.Range("D61").Value = 44771
$ws.Range("J61").Value = 35
$ws.Range("K61").Value = 40000
$ws.Range("L61").Value = 40000
$ws.Range("M61").Value = 40000
$ws.Range("N61").Value = '$/malla 25 kilos'
$ws.Range("O61").Value = 'Perú'
$ws.Range("P61").Value = 1600
$ws.Range("D62").Value = 44747
$ws.Range("K62").Value = 27000
$ws.Range("L62").Value = 27000
$ws.Range("M62").Value = 27000
$ws.Range("P62").Value = 1080
$ws.Range("D63").Value = 44855
$ws.Range("K63").Value = 35000
$ws.Range("L63").Value = 35000
$ws.Range("M63").Value = 35000
$ws.Range("P63").Value = 1400
$ws.Range("D64").Value = 45079
$ws.Range("J64").Value = 40
$ws.Range("K64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("M64").Value = 30000
$ws.Range("P64").Value = 1200
